$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1449608478139278
$ws.Range("C2").Value = 0.3375693222279904
$ws.Range("D2").Value = 0.2207894770556722
$ws.Range("E2").Value = 0.4698824076890645
$ws.Range("F2").Value = 0.4638353168623497
$ws.Range("B3").Value = 0.1842546231075045
$ws.Range("C3").Value = 0.2869443818724057
$ws.Range("D3").Value = 0.1607157978134096
$ws.Range("E3").Value = 0.4008937487831527
$ws.Range("F3").Value = 0.3705804469342626
$ws.Range("B4").Value = 0.2080690016882378
$ws.Range("C4").Value = 0.2875078072560339
$ws.Range("D4").Value = 0.1982701821409898
$ws.Range("E4").Value = 0.4452754003321874
$ws.Range("F4").Value = 0.4111767671329999
$ws.Range("B5").Value = 0.200282133760416
$ws.Range("C5").Value = 0.2861112695554327
$ws.Range("D5").Value = 0.1811512371320833
$ws.Range("E5").Value = 0.4256186522370504
$ws.Range("F5").Value = 0.3938808632458601
$ws.Range("B6").Value = 0.2018315634777409
$ws.Range("C6").Value = 0.2396510267933027
$ws.Range("D6").Value = 0.1874691252642877
$ws.Range("E6").Value = 0.4329770493505258
$ws.Range("F6").Value = 0.4037781916520482
$ws.Range("B7").Value = 0.2248954171260089
$ws.Range("C7").Value = 0.2511991784966441
$ws.Range("D7").Value = 0.1348234142648089
$ws.Range("E7").Value = 0.3671830800361162
$ws.Range("F7").Value = 0.3078573514196035
$ws.Range("B8").Value = 0.2892781581597261
$ws.Range("C8").Value = 0.3217998310022198
$ws.Range("D8").Value = 0.2884461656711942
$ws.Range("E8").Value = 0.5370718440499317
$ws.Range("F8").Value = 0.4956986740546043
$ws.Range("B9").Value = 0.527914141851893
$ws.Range("C9").Value = 0.527914141851893
$ws.Range("D9").Value = 0.3863563944451271
$ws.Range("E9").Value = 0.6215757350839294
$ws.Range("F9").Value = 0.4018638823244256
$ws.Range("B10").Value = 0.01851010070983469
$ws.Range("C10").Value = 0.01851010070983469
$ws.Range("D10").Value = 0.0003426238282882228
$ws.Range("E10").Value = 0.01851010070983469
